$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.645.95"
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = "'2.282.40"
$ws.Range('E3').Value = '  -3.83%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'300.67"
$ws.Range('E5').Value = '  -3.00%  '
$ws.Range('D6').Value = "'97.41"
$ws.Range('E6').Value = '  -6.40%  '
$ws.Range('E7').Value = '  -1.93%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = "'0.499"
$ws.Range('E9').Value = '  -4.11%  '
$ws.Range('D10').Value = "'33.62"
$ws.Range('E10').Value = '  -5.95%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = "'50.84"
$ws.Range('E11').Value = '  -4.51%  '
$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = "'0.0789"
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('E14').Value = '  -4.09%  '
$ws.Range('D15').Value = "'2.636.47"
$ws.Range('E15').Value = '  -3.90%  '
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').Value = "'2.266.00"
$ws.Range('E17').Value = '  -4.54%  '
$ws.Range('D18').Value = "'0.789"
$ws.Range('E18').Value = '  -2.77%  '
$ws.Range('D19').Value = "'42.515.73"
$ws.Range('E19').Value = '  -2.15%  '
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').Value = "'11.45"
$ws.Range('E21').Value = '  -3.80%  '
$ws.Range('D22').Value = "'6.00"
$ws.Range('E22').Value = '  -4.99%  '
$ws.Range('D23').Value = "'66.69"
$ws.Range('E23').Value = '  -2.23%  '
$ws.Range('D24').Value = "'235.25"
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('D25').Value = "'1.94"
$ws.Range('E25').Value = '  -5.07%  '
$ws.Range('D26').Value = "'2.49"
$ws.Range('E26').Value = '  -4.26%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = "'24.43"
$ws.Range('E28').Value = '  -4.77%  '
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D30').Value = "'164.68"
$ws.Range('E30').Value = '  +2.25%  '
$ws.Range('D31').Value = "'33.62"
$ws.Range('E31').Value = '  -7.90%  '
$ws.Range('D32').Value = "'9.11"
$ws.Range('E32').Value = '  -3.83%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = "'4.96"
$ws.Range('E34').Value = '  -4.84%  '
$ws.Range('E35').Value = '  -4.52%  '
$ws.Range('D37').Value = "'4.34"
$ws.Range('E37').Value = '  -6.81%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = "'16.24"
$ws.Range('E38').Value = '  -10.70%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = "'2.82"
$ws.Range('E39').Value = '  -8.19%  '
$ws.Range('D40').Value = "'0.100"
$ws.Range('E40').Value = '  -4.90%  '
$ws.Range('E41').Value = '  -8.02%  '
$ws.Range('E42').Value = '  -3.24%  '
$ws.Range('E43').Value = '  -7.98%  '
$ws.Range('D44').Value = "'1.959.63"
$ws.Range('E44').Value = '  -3.56%  '
$ws.Range('E45').Value = '  -2.29%  '
$ws.Range('D46').Value = "'17.78"
$ws.Range('E46').Value = '  -9.44%  '
$ws.Range('D47').Value = "'9.69"
$ws.Range('E47').Value = '  -8.03%  '
$ws.Range('E48').Value = '  -8.58%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = "'53.34"
$ws.Range('E49').Value = '  -7.52%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = "'2.82"
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = "'2.504.51"
$ws.Range('E51').Value = '  -4.06%  '
